$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet is protected; temporarily unprotect so the cells below can be edited.
$ws.Unprotect()

# Update the confidential-notice date text (A10): 2021-05-05 -> 2021-05-06
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-06 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-7
$ws.Range("D2").Value = 0.2459807335371519
$ws.Range("E2").Value = 0.001836087464530234

$ws.Range("D3").Value = 0.4978864535132011
$ws.Range("E3").Value = 0.008438818565400963

$ws.Range("D4").Value = 0.09739919023087437
$ws.Range("E4").Value = -0.01696898771211242

$ws.Range("D5").Value = 0.1019665747128018
$ws.Range("E5").Value = 0.007598371777476354

$ws.Range("D6").Value = 0.05676704800597086
$ws.Range("E6").Value = 0.004060455673358998

$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0.004005729952214265

# Restore sheet protection (matches original protected state).
$ws.Protect()
